$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Sample_0100", 0),
    @("Sample_0226", 1),
    @("Sample_0085", 1),
    @("Sample_0047", 1),
    @("Sample_0154", 0),
    @("Sample_0218", 1),
    @("Sample_0088", 1),
    @("Sample_0251", 1),
    @("Sample_0160", 1),
    @("Sample_0038", 0),
    @("Sample_0206", 0),
    @("Sample_0208", 0),
    @("Sample_0242", 1),
    @("Sample_0006", 0),
    @("Sample_0093", 1),
    @("Sample_0163", 1),
    @("Sample_0204", 0),
    @("Sample_0219", 1),
    @("Sample_0074", 0),
    @("Sample_0246", 1),
    @("Sample_0033", 1),
    @("Sample_0187", 1),
    @("Sample_0148", 1),
    @("Sample_0053", 1),
    @("Sample_0018", 0),
    @("Sample_0185", 1),
    @("Sample_0215", 1),
    @("Sample_0013", 1),
    @("Sample_0039", 1),
    @("Sample_0168", 0),
    @("Sample_0011", 1),
    @("Sample_0052", 1),
    @("Sample_0253", 1),
    @("Sample_0201", 0),
    @("Sample_0157", 1),
    @("Sample_0200", 0),
    @("Sample_0017", 0),
    @("Sample_0115", 0),
    @("Sample_0159", 1),
    @("Sample_0202", 0),
    @("Sample_0003", 1),
    @("Sample_0235", 1),
    @("Sample_0238", 1),
    @("Sample_0124", 0),
    @("Sample_0170", 0),
    @("Sample_0117", 1),
    @("Sample_0091", 1),
    @("Sample_0186", 1),
    @("Sample_0040", 0),
    @("Sample_0087", 1),
    @("Sample_0176", 0),
    @("Sample_0240", 1),
    @("Sample_0212", 1),
    @("Sample_0037", 1),
    @("Sample_0034", 1),
    @("Sample_0207", 0),
    @("Sample_0203", 0),
    @("Sample_0139", 1),
    @("Sample_0167", 0),
    @("Sample_0137", 0),
    @("Sample_0142", 0),
    @("Sample_0086", 1),
    @("Sample_0112", 1),
    @("Sample_0180", 1),
    @("Sample_0144", 1),
    @("Sample_0135", 0),
    @("Sample_0250", 1),
    @("Sample_0224", 1),
    @("Sample_0094", 1),
    @("Sample_0152", 0),
    @("Sample_0210", 0),
    @("Sample_0090", 0),
    @("Sample_0234", 1),
    @("Sample_0009", 1),
    @("Sample_0164", 1),
    @("Sample_0103", 0),
    @("Sample_0193", 1),
    @("Sample_0236", 1),
    @("Sample_0213", 1),
    @("Sample_0241", 1),
    @("Sample_0097", 1),
    @("Sample_0138", 1),
    @("Sample_0057", 1),
    @("Sample_0014", 0),
    @("Sample_0169", 0),
    @("Sample_0149", 0),
    @("Sample_0141", 0),
    @("Sample_0232", 1),
    @("Sample_0237", 1),
    @("Sample_0147", 1),
    @("Sample_0015", 1),
    @("Sample_0020", 1),
    @("Sample_0110", 0),
    @("Sample_0227", 1),
    @("Sample_0191", 1),
    @("Sample_0025", 1),
    @("Sample_0126", 0),
    @("Sample_0231", 1),
    @("Sample_0194", 1),
    @("Sample_0223", 0),
    @("Sample_0254", 1),
    @("Sample_0116", 0),
    @("Sample_0105", 1),
    @("Sample_0166", 0),
    @("Sample_0099", 0),
    @("Sample_0252", 1),
    @("Sample_0109", 1),
    @("Sample_0095", 0),
    @("Sample_0146", 1),
    @("Sample_0032", 1),
    @("Sample_0108", 1),
    @("Sample_0220", 1),
    @("Sample_0249", 1),
    @("Sample_0189", 1),
    @("Sample_0083", 0),
    @("Sample_0217", 1),
    @("Sample_0244", 1),
    @("Sample_0004", 1),
    @("Sample_0101", 0),
    @("Sample_0209", 0),
    @("Sample_0190", 1),
    @("Sample_0211", 0),
    @("Sample_0162", 1),
    @("Sample_0158", 1),
    @("Sample_0175", 0),
    @("Sample_0145", 0),
    @("Sample_0131", 1),
    @("Sample_0181", 1),
    @("Sample_0255", 1),
    @("Sample_0031", 1),
    @("Sample_0248", 1),
    @("Sample_0125", 0),
    @("Sample_0165", 0),
    @("Sample_0022", 1),
    @("Sample_0150", 0),
    @("Sample_0239", 1),
    @("Sample_0161", 1)
)

$startRow = 2
$oldLastRow = 134

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    if ($row -gt $oldLastRow) {
        $ws.Range("A$oldLastRow").Copy() | Out-Null
        $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    }
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
